# Update the cryptocurrency price/volume table (rows 2-51) to reflect the
# latest scrape, matching the GitHub Actions "Updated cryptos list" commit.
# Also swaps the Uniswap/WrappedEther rows (20 and 21), which traded ranking
# positions in this run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the data range to plain text storage before writing so numeric-looking
# strings (e.g. "0.9996", "27.226.56") are not auto-coerced into Excel numbers,
# matching the original inline-string cell type. Style is reset back to Normal
# afterwards so no stray number-format/style is left behind on the cells.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.226.56"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "1.903.80"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "306.02"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").Value = "0.5385"
$ws.Range("E7").Value = "  +3.15%  "
$ws.Range("D8").Value = "0.3807"
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("D10").Value = "22.16"
$ws.Range("E10").Value = "  +4.98%  "
$ws.Range("D11").Value = "0.9050"
$ws.Range("E11").Value = "  +0.57%  "
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("D13").Value = "95.74"
$ws.Range("E13").Value = "  -0.63%  "
$ws.Range("D14").Value = "5.348"
$ws.Range("E14").Value = "  +1.18%  "
$ws.Range("D15").Value = "0.9990"
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("E16").Value = "  +2.11%  "
$ws.Range("D17").Value = "0.000008667"
$ws.Range("E17").Value = "  +1.07%  "
$ws.Range("D18").Value = "0.9993"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").Value = "27.261.43"
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "5.047"
$ws.Range("E20").Value = "  -0.65%  "
$ws.Range("B21").Value = "WrappedEther"
$ws.Range("C21").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D21").Value = "1.118.43"
$ws.Range("E21").Value = "  -41.89%  "
$ws.Range("E22").Value = "  +0.99%  "
$ws.Range("D23").Value = "6.517"
$ws.Range("E23").Value = "  +1.81%  "
$ws.Range("D24").Value = "148.44"
$ws.Range("E24").Value = "  +0.54%  "
$ws.Range("D25").Value = "2.305"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("D26").Value = "18.39"
$ws.Range("E26").Value = "  +1.20%  "
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("D28").Value = "116.74"
$ws.Range("E28").Value = "  +1.57%  "
$ws.Range("D29").Value = "4.852"
$ws.Range("E29").Value = "  +1.34%  "
$ws.Range("D30").Value = "4.701"
$ws.Range("E30").Value = "  -3.91%  "
$ws.Range("D31").Value = "0.09217"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").Value = "0.8272"
$ws.Range("E32").Value = "  +4.92%  "
$ws.Range("D33").Value = "0.05079"
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("D34").Value = "1.218"
$ws.Range("E34").Value = "  +0.27%  "
$ws.Range("D35").Value = "3.003"
$ws.Range("E35").Value = "  +1.51%  "
$ws.Range("D36").Value = "3.319"
$ws.Range("E36").Value = "  -3.01%  "
$ws.Range("D37").Value = "2.682"
$ws.Range("E37").Value = "  +3.89%  "
$ws.Range("D38").Value = "0.5916"
$ws.Range("E38").Value = "  +4.08%  "
$ws.Range("D39").Value = "0.01999"
$ws.Range("E39").Value = "  +0.56%  "
$ws.Range("E40").Value = "  +0.57%  "
$ws.Range("D41").Value = "9.275"
$ws.Range("E41").Value = "  +2.88%  "
$ws.Range("D42").Value = "6.666"
$ws.Range("E42").Value = "  +1.69%  "
$ws.Range("D43").Value = "116.65"
$ws.Range("E43").Value = "  +0.45%  "
$ws.Range("D44").Value = "0.5109"
$ws.Range("E44").Value = "  +5.03%  "
$ws.Range("D45").Value = "0.1528"
$ws.Range("E45").Value = "  +0.68%  "
$ws.Range("D46").Value = "10.17"
$ws.Range("E46").Value = "  +1.11%  "
$ws.Range("D47").Value = "0.9990"
$ws.Range("E47").Value = "  -0.25%  "
$ws.Range("E48").Value = "  +1.20%  "
$ws.Range("D49").Value = "38.26"
$ws.Range("E49").Value = "  +0.43%  "
$ws.Range("D50").Value = "0.06108"
$ws.Range("E50").Value = "  +2.99%  "
$ws.Range("D51").Value = "63.42"
$ws.Range("E51").Value = "  +0.06%  "

$dataRange.Style = "Normal"
